# This script updates the "想去人数" (number of people interested) figures
# in column F on the "展览" sheet and the "全部类型" sheet, reflecting a
# refreshed data scrape (per commit message: "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F4").Value = 94
$ws1.Range("F5").Value = 44
$ws1.Range("F8").Value = 36
$ws1.Range("F9").Value = 8277
$ws1.Range("F10").Value = 773
$ws1.Range("F11").Value = 278
$ws1.Range("F12").Value = 1118
$ws1.Range("F13").Value = 855
$ws1.Range("F14").Value = 55
$ws1.Range("F15").Value = 38
$ws1.Range("F16").Value = 212
$ws1.Range("F17").Value = 114
$ws1.Range("F18").Value = 55
$ws1.Range("F20").Value = 897

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F4").Value = 94
$ws4.Range("F6").Value = 44
$ws4.Range("F10").Value = 36
$ws4.Range("F11").Value = 8278
$ws4.Range("F12").Value = 773
$ws4.Range("F13").Value = 278
$ws4.Range("F14").Value = 1118
$ws4.Range("F15").Value = 855
$ws4.Range("F16").Value = 55
$ws4.Range("F17").Value = 38
$ws4.Range("F18").Value = 212
$ws4.Range("F19").Value = 114
$ws4.Range("F20").Value = 55
$ws4.Range("F22").Value = 897
